# Updates cryptos list data (prices / 1h volume %) and reorders a few rows
# to match the refreshed scrape, per commit "Updated cryptos list ... with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value-only updates (rows whose coin identity/order is unchanged) ---

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.302.95"
$ws.Range("E2").Value = "  -0.20%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.482.58"
$ws.Range("E3").Value = "  -1.07%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "568.03"
$ws.Range("E5").Value = "  -0.83%  "

# Row 6 - Solana
$ws.Range("D6").Value = "163.27"
$ws.Range("E6").Value = "  -1.83%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.04%  "

# Row 8 - XRP (force text so the trailing zero in "0.510" survives)
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.510"
$ws.Range("E8").Value = "  -1.10%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.481.86"
$ws.Range("E9").Value = "  -1.05%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.157"
$ws.Range("E10").Value = "  -1.95%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.70%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "0.352"
$ws.Range("E12").Value = "  -1.31%  "

# Row 13 - Toncoin
$ws.Range("D13").Value = "4.87"
$ws.Range("E13").Value = "  -1.26%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.939.10"
$ws.Range("E14").Value = "  -1.31%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "69.179.81"
$ws.Range("E15").Value = "  -0.48%  "

# Row 16 - ShibaInu
$ws.Range("D16").Value = "0.0000173"
$ws.Range("E16").Value = "  -1.19%  "

# Row 17 - Avalanche
$ws.Range("D17").Value = "24.02"
$ws.Range("E17").Value = "  -3.36%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.514.96"
$ws.Range("E18").Value = "  +0.22%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  -2.04%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "7.33"
$ws.Range("E20").Value = "  -4.46%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "345.36"
$ws.Range("E21").Value = "  -1.11%  "

# Row 22 - Polkadot
$ws.Range("D22").Value = "3.87"
$ws.Range("E22").Value = "  -1.53%  "

# Row 23 - SuiNetwork
$ws.Range("D23").Value = "1.88"
$ws.Range("E23").Value = "  -5.51%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.04%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "69.32"
$ws.Range("E25").Value = "  -0.98%  "

# Row 26 - NEARProtocol
$ws.Range("D26").Value = "3.86"
$ws.Range("E26").Value = "  -2.94%  "

# Row 27 - WrappedeETH : unchanged

# Row 28 - Aptos
$ws.Range("D28").Value = "8.56"
$ws.Range("E28").Value = "  -3.85%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  -0.03%  "

# Row 30 - PEPE (contains subscript-3 char U+2083)
$ws.Range("D30").Value = "0.0$([char]0x2083)0860"
$ws.Range("E30").Value = "  -3.96%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "7.52"
$ws.Range("E31").Value = "  -4.67%  "

# --- Rows 32/33 swap: Bittensor <-> Fetch.AI (order flips, data refreshed) ---
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.18"
$ws.Range("E32").Value = "  -4.88%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "434.62"
$ws.Range("E33").Value = "  -6.22%  "

# Row 34 - FirstDigitalUSD
$ws.Range("E34").Value = "  -0.03%  "

# Row 35 - PancakeSwap
$ws.Range("E35").Value = "  -2.08%  "

# Row 36 - Monero
$ws.Range("D36").Value = "156.83"
$ws.Range("E36").Value = "  -0.41%  "

# --- Rows 37/38 swap: Kaspa <-> WhiteBITCoin (order flips, data refreshed) ---
$ws.Range("B37").Value = "WhiteBITCoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D37").Value = "19.05"
$ws.Range("E37").Value = "  +0.16%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.112"
$ws.Range("E38").Value = "  -3.44%  "

# Row 39 - EthereumClassic
$ws.Range("D39").Value = "18.04"
$ws.Range("E39").Value = "  -2.40%  "

# Row 40 - USDe
$ws.Range("E40").Value = "  +0.02%  "

# Row 41 - PolygonEcosystemToken
$ws.Range("D41").Value = "0.311"
$ws.Range("E41").Value = "  -2.07%  "

# Row 42 - RenderToken
$ws.Range("D42").Value = "4.54"
$ws.Range("E42").Value = "  -4.59%  "

# Row 43 - Stacks
$ws.Range("D43").Value = "1.56"
$ws.Range("E43").Value = "  -2.65%  "

# Row 44 - dogwifhat
$ws.Range("D44").Value = "2.14"
$ws.Range("E44").Value = "  -6.07%  "

# Row 45 - ImmutableX
$ws.Range("E45").Value = "  -6.42%  "

# Row 46 - Aave
$ws.Range("D46").Value = "137.37"
$ws.Range("E46").Value = "  -3.01%  "

# Row 47 - Filecoin (force text so the trailing zero in "3.40" survives)
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.40"
$ws.Range("E47").Value = "  -2.39%  "

# --- Rows 48/49 swap: ARBITRUM <-> POPCAT (order flips, data refreshed) ---
$ws.Range("B48").Value = "POPCAT"
$ws.Range("C48").Value = "https://coinranking.com/coin/sLBuDEsp6+popcat-popcat"
$ws.Range("D48").Value = "2.06"
$ws.Range("E48").Value = "  +27.51%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "0.502"
$ws.Range("E49").Value = "  -4.49%  "

# Row 50 - Cronos
$ws.Range("D50").Value = "0.0722"
$ws.Range("E50").Value = "  -1.07%  "

# Row 51 - Mantle (force text so the trailing zero in "0.570" survives)
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.570"
$ws.Range("E51").Value = "  -0.97%  "
